$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.422.76"

$ws.Range("D3").Value = "2.981.62"
$ws.Range("E3").Value = "  -2.64%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.12"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.69"
$ws.Range("E6").Value = "  +4.89%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  -0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.31"
$ws.Range("E9").Value = "  +3.52%  "

$ws.Range("E10").Value = "  +2.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.350"
$ws.Range("E11").Value = "  -2.47%  "

$ws.Range("E12").Value = "  +0.25%  "

$ws.Range("D13").Value = "3.490.93"
$ws.Range("E13").Value = "  -3.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.10"
$ws.Range("E14").Value = "  +3.84%  "

$ws.Range("D15").Value = "56.405.20"
$ws.Range("E15").Value = "  +3.17%  "

$ws.Range("E16").Value = "  +3.86%  "

$ws.Range("D17").Value = "2.981.22"
$ws.Range("E17").Value = "  -3.51%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.68"
$ws.Range("E18").Value = "  +3.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.34"
$ws.Range("E19").Value = "  -0.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.76"
$ws.Range("E20").Value = "  +2.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "325.30"
$ws.Range("E21").Value = "  -0.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("E23").Value = "  -3.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.16"
$ws.Range("E24").Value = "  -5.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.12%  "

$ws.Range("E26").Value = "  -1.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("E29").Value = "  -0.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.77"
$ws.Range("E30").Value = "  +2.88%  "

$ws.Range("E31").Value = "  -3.46%  "

$ws.Range("E32").Value = "  -2.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.14"
$ws.Range("E34").Value = "  +0.00%  "

$ws.Range("E35").Value = "  -2.48%  "

$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("E37").Value = "  -5.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0672"
$ws.Range("E38").Value = "  +2.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "23.01"
$ws.Range("E39").Value = "  +0.67%  "

$ws.Range("D40").Value = "3.014.33"
$ws.Range("E40").Value = "  -3.15%  "

$ws.Range("E41").Value = "  -0.12%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.05"
$ws.Range("E42").Value = "  -0.38%  "

$ws.Range("E43").Value = "  -3.28%  "

$ws.Range("D44").Value = "2.235.20"
$ws.Range("E44").Value = "  +0.74%  "

$ws.Range("E45").Value = "  -4.72%  "

$ws.Range("E46").Value = "  +1.79%  "

$ws.Range("E47").Value = "  -3.98%  "

$ws.Range("E48").Value = "  +12.95%  "

$ws.Range("E49").Value = "  +3.98%  "

$ws.Range("E50").Value = "  -2.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.90"
$ws.Range("E51").Value = "  -3.11%  "
